# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking "Price" values are entered with a leading apostrophe and the
# cell style is reset to "Normal" afterwards so they stay plain text cells
# (matching the sheet's existing inline/shared-string layout) instead of being
# auto-coerced to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.424.44"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.573.39"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'212.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'44.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("D9").Value = "'23.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'0.0894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "1.798.09"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "1.565.02"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "28.417.33"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "'0.516"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'230.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'3.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "'151.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'14.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'0.0484"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "1.381.24"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'0.521"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "'5.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("D47").Value = "'62.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'0.919"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.21%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.710.40"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "'85.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.98%  "
